# Generate Report for Handoff
#
# The localization-status report has moved from "In Translation" to
# "Ready for handoff": update the Status cells (and the Overview rollup
# columns that mirror them) plus the "Latest Handoff Datetime" timestamps
# that were stamped when the handoff package was produced. Then refresh
# the column widths that depend on the new "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-26 20:57:08"

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-26 20:56:58"

# ---- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-26 20:57:08"

# ---- Resize the Status columns so the longer "Ready for handoff" text
#      is not truncated (mirrors the width refresh captured in the diff).
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
